# compressor compressibility factor add
# Sheet1 layout: column B = isen_eff, rows -> A2 piston, A3 diaphragm, A4 centrifugal
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the previously-blank isen_eff (compressibility) factor for "piston" (row 2)
$ws.Range("B2").Value = 0.82

# Update the isen_eff (compressibility) factor for "centrifugal" (row 4)
$ws.Range("B4").Value = 0.78

# Leave the active selection on B4, the last cell edited
$ws.Range("B4").Select()
